$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4783.16
$ws.Range("C3").Value = 4860.45
$ws.Range("C4").Value = 4939.63
$ws.Range("C5").Value = 5020.68
$ws.Range("C6").Value = 5101.72
$ws.Range("C7").Value = 5180.74
$ws.Range("C8").Value = 5261.27
